# Region3_Stations_V1 - "adjusted Atl City and Cape May node indices"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Region3_Stations_V1")
$ws.Activate()

# --- Row 16 (Atlantic City, NJ): updated nearest Region3 node lon/lat and node indices ---
$ws.Range("K16").Value = -74.418297269295806
$ws.Range("L16").Value = 39.3567908537158
$ws.Range("N16").Value = 1942291
$ws.Range("O16").Value = 964479

# --- Row 17 (Cape May Ferry Terminal, NJ): updated nearest Region3 node lon/lat/index + row no longer highlighted ---
$ws.Range("K17").Value = -74.959962000000004
$ws.Range("L17").Value = 38.967747000000003
$ws.Range("N17").Value = 144964
$ws.Range("O17").Value = 72477

# remove the highlight fill from row 17 (match the unhighlighted style used by nearby rows)
$ws.Range("A18:J18").Copy()
$ws.Range("A17:J17").PasteSpecial(-4122)
$ws.Range("M18:O18").Copy()
$ws.Range("M17:O17").PasteSpecial(-4122)
$ws.Range("K18:L18").Copy()
$ws.Range("K17:L17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 25 (Springmaid Pier, SC): remove highlight fill + set explicit row height (values unchanged) ---
$ws.Range("A26:J26").Copy()
$ws.Range("A25:J25").PasteSpecial(-4122)
$ws.Range("M26:O26").Copy()
$ws.Range("M25:O25").PasteSpecial(-4122)
$ws.Range("K26:L26").Copy()
$ws.Range("K25:L25").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Rows.Item(25).RowHeight = 17

# --- Selection moved to O16 ---
$ws.Range("O16").Select()
